$d = $word.ActiveDocument

# The edit merges paragraphs of the form "<id>XXXX</id>" (currently split across
# three runs: "<id>", "XXXX", "</id>") into a single run containing the full
# "<id>XXXX</id>" text (keeping the first run's formatting), while leaving the
# trailing empty run untouched. This pattern occurs for each page/item id
# marker in the document (p110r_3, p110v_1, p110v_2, p110v_3, ...).

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)

    if ($t.StartsWith("<id>") -and $t.EndsWith("</id>")) {
        $rng = $p.Range
        $rng.Find.Execute($t, $true, $false, $false, $false, $false, $true, 1, $false, $t, 2)
    }
}
